$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# New mod setups: bump the base deposit tier (C2) and switch on the H-column
# multiplier (H2) — this cascades through all the dependent formulas below.
$ws.Range("C2").Value = 5
$ws.Range("H2").Value = 1

# Update the saved view state (scroll position + active selection) to match.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 129
$ws.Range("H3").Select()
